$wb = $excel.ActiveWorkbook

# Sheet R1 (sheet1.xml)
$ws = $wb.Worksheets.Item("R1")
$ws.Range("G2").Value = "3926:12:40"
$ws.Range("G3").Value = "65:45:18"

# Sheet R2 (sheet2.xml)
$ws = $wb.Worksheets.Item("R2")
$ws.Range("G2").Value = "12107:36:19"
$ws.Range("G3").Value = "3237:19:48"
$ws.Range("G4").Value = "475:31:22"

# Sheet R4 (sheet4.xml)
$ws = $wb.Worksheets.Item("R4")
$ws.Range("G2").Value = "2953:26:08"
$ws.Range("G3").Value = "180:38:23"

# Sheet R5 (sheet5.xml)
$ws = $wb.Worksheets.Item("R5")
$ws.Range("G2").Value = "427:25:07"

# Sheet R6 (sheet6.xml)
$ws = $wb.Worksheets.Item("R6")
$ws.Range("G2").Value = "67:57:25"
